$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: rename measurement point ---
$ws.Range("A3").Value = "Camera di manovra Merone"

# --- Row 2 & 3: corrected lat/long positions ---
$ws.Range("F2").Value = 39.860277777777782
$ws.Range("G2").Value = 16.508611111111112
$ws.Range("F3").Value = 39.176666666666669
$ws.Range("G3").Value = 16.341666666666665

# --- Row 9: status change ---
$ws.Range("K9").Value = "Inacessibile"

# --- Rows 10-13: corrected names / positions for "In valutazione" measurers ---
$ws.Range("F10").Value = 39.240833333333335
$ws.Range("G10").Value = 16.361944444444447

$ws.Range("A11").Value = "Sorgente Zumpo"
$ws.Range("F11").Value = 39.231111111111112
$ws.Range("G11").Value = 16.404999999999998

$ws.Range("A12").Value = "SA"
$ws.Range("F12").Value = 38.999938888888892
$ws.Range("G12").Value = 17.062283333333333

$ws.Range("A13").Value = "Sorical Differenzi Murate / Pisarello"
$ws.Range("F13").Value = 39.116572222222224
$ws.Range("G13").Value = 16.749461111111113

# --- Move the bottom spacer row from 23 to 22 (one new data row inserted) ---
$ws.Rows.Item(23).Delete()

# --- Row 14: new measurer in evaluation ---
$ws.Range("A14").Value = "I salto Acquedotto Merone"
$ws.Range("F14").Formula = "=39+10/60+55/3600"
$ws.Range("G14").Formula = "=16+20/60+9/3600"
$ws.Range("J14").Value = "areatecnica"
$ws.Range("K14").Value = "In valutazione"

$ws.Rows.Item(22).RowHeight = 9

# --- Column K width (matches new longer content) ---
$ws.Columns.Item(11).AutoFit()

# --- sheet view tweaks ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.Zoom = 145
$ws.Range("L10").Select()
